# Auto-generated edit script: updates market-price-derived columns (H-N)
# across multiple item-crafting sheets, per the commit's scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 299.13333
$ws.Range("I28").Value = 187.92308
$ws.Range("K28").Value = 187.92308
$ws.Range("M28").Value = 297.07692
$ws.Range("H32").Value = 466.66666
$ws.Range("J32").Value = 466.66666
$ws.Range("L32").Value = 466.66666
$ws.Range("N32").Value = -1118.66666
$ws.Range("H86").Value = 6680.3335
$ws.Range("I86").Value = 1279
$ws.Range("J86").Value = 20723.8
$ws.Range("K86").Value = 1279
$ws.Range("L86").Value = 20723.8
$ws.Range("M86").Value = -156
$ws.Range("N86").Value = -22969.8
$ws.Range("H89").Value = 6680.3335
$ws.Range("I89").Value = 1279
$ws.Range("J89").Value = 20723.8
$ws.Range("K89").Value = 6395
$ws.Range("L89").Value = 103619
$ws.Range("M89").Value = -779
$ws.Range("N89").Value = -114851
$ws.Range("H106").Value = 11907974
$ws.Range("I106").Value = 47621020
$ws.Range("J106").Value = 3624.6667
$ws.Range("K106").Value = 47621020
$ws.Range("L106").Value = 3624.6667
$ws.Range("M106").Value = -47620389
$ws.Range("N106").Value = -4886.6667
$ws.Range("H116").Value = 4264.5713
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 4400.3076
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 4400.3076
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -11284.3076
$ws.Range("H137").Value = 1567.6428
$ws.Range("I137").Value = 1612.2
$ws.Range("K137").Value = 4836.6
$ws.Range("M137").Value = -2286.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8042.068
$ws.Range("I32").Value = 6878.439
$ws.Range("J32").Value = 10692.556
$ws.Range("K32").Value = 6878.439
$ws.Range("L32").Value = 10692.556
$ws.Range("M32").Value = -6591.439
$ws.Range("N32").Value = -11266.556
$ws.Range("H61").Value = 3128.2258
$ws.Range("I61").Value = 2891.3462
$ws.Range("K61").Value = 2891.3462
$ws.Range("M61").Value = -2679.3462
$ws.Range("H63").Value = 1646165
$ws.Range("I63").Value = 1537.2941
$ws.Range("J63").Value = 15625500
$ws.Range("K63").Value = 1537.2941
$ws.Range("L63").Value = 15625500
$ws.Range("M63").Value = -851.2941000000001
$ws.Range("N63").Value = -15626872
$ws.Range("H66").Value = 1646165
$ws.Range("I66").Value = 1537.2941
$ws.Range("J66").Value = 15625500
$ws.Range("K66").Value = 7686.4705
$ws.Range("L66").Value = 78127500
$ws.Range("M66").Value = -4254.4705
$ws.Range("N66").Value = -78134364
$ws.Range("H124").Value = 9224.5
$ws.Range("J124").Value = 9224.5
$ws.Range("L124").Value = 9224.5
$ws.Range("N124").Value = -19044.5
$ws.Range("H125").Value = 31992.666
$ws.Range("J125").Value = 31992.666
$ws.Range("L125").Value = 31992.666
$ws.Range("N125").Value = -41832.666
$ws.Range("H132").Value = 17108.97
$ws.Range("I132").Value = 2365.25
$ws.Range("J132").Value = 38171.43
$ws.Range("K132").Value = 7095.75
$ws.Range("L132").Value = 114514.29
$ws.Range("M132").Value = -4565.75
$ws.Range("N132").Value = -119574.29
$ws.Range("H135").Value = 32405.2
$ws.Range("J135").Value = 32405.2
$ws.Range("L135").Value = 32405.2
$ws.Range("N135").Value = -42545.2
$ws.Range("H136").Value = 3128.2258
$ws.Range("I136").Value = 2891.3462
$ws.Range("K136").Value = 8674.0386
$ws.Range("M136").Value = -6124.0386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 673.6667
$ws.Range("I22").Value = 611.36365
$ws.Range("J22").Value = 845
$ws.Range("K22").Value = 611.36365
$ws.Range("L22").Value = 845
$ws.Range("M22").Value = -438.36365
$ws.Range("N22").Value = -1191

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3751.7097
$ws.Range("I31").Value = 939.2308
$ws.Range("J31").Value = 5782.9443
$ws.Range("K31").Value = 939.2308
$ws.Range("L31").Value = 5782.9443
$ws.Range("M31").Value = -644.2308
$ws.Range("N31").Value = -6372.9443
$ws.Range("H34").Value = 3751.7097
$ws.Range("I34").Value = 939.2308
$ws.Range("J34").Value = 5782.9443
$ws.Range("K34").Value = 939.2308
$ws.Range("L34").Value = 5782.9443
$ws.Range("M34").Value = -737.2308
$ws.Range("N34").Value = -6186.9443
$ws.Range("H105").Value = 1189.3684
$ws.Range("I105").Value = 969.13336
$ws.Range("K105").Value = 969.13336
$ws.Range("M105").Value = 777.86664
$ws.Range("H107").Value = 1459.8636
$ws.Range("I107").Value = 884.7273
$ws.Range("J107").Value = 2035
$ws.Range("K107").Value = 884.7273
$ws.Range("L107").Value = 2035
$ws.Range("M107").Value = 1035.2727
$ws.Range("N107").Value = -5875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1294.3871
$ws.Range("I5").Value = 1127.7142
$ws.Range("K5").Value = 3383.1426
$ws.Range("M5").Value = -3271.1426
$ws.Range("H12").Value = 124.28571
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 166.2
$ws.Range("K12").Value = 58.5
$ws.Range("L12").Value = 498.6
$ws.Range("M12").Value = 114.5
$ws.Range("N12").Value = -844.5999999999999
$ws.Range("H45").Value = 510
$ws.Range("I45").Value = 265
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 795
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -263
$ws.Range("N45").Value = -4064
$ws.Range("H122").Value = 1173.32
$ws.Range("I122").Value = 336
$ws.Range("J122").Value = 1332.8096
$ws.Range("K122").Value = 3024
$ws.Range("L122").Value = 11995.2864
$ws.Range("M122").Value = -574
$ws.Range("N122").Value = -16895.2864
$ws.Range("H131").Value = 736.0303
$ws.Range("I131").Value = 576.6667
$ws.Range("J131").Value = 741.01044
$ws.Range("K131").Value = 1730.0001
$ws.Range("L131").Value = 2223.03132
$ws.Range("M131").Value = 3309.9999
$ws.Range("N131").Value = -12303.03132
$ws.Range("H135").Value = 1294.3871
$ws.Range("I135").Value = 1127.7142
$ws.Range("K135").Value = 10149.4278
$ws.Range("M135").Value = -7614.427799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3481294.8
$ws.Range("I70").Value = 3950.8
$ws.Range("K70").Value = 3950.8
$ws.Range("M70").Value = -3680.8
$ws.Range("H73").Value = 3481294.8
$ws.Range("I73").Value = 3950.8
$ws.Range("K73").Value = 3950.8
$ws.Range("M73").Value = -3014.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1157.48
$ws.Range("I113").Value = 1266.7142
$ws.Range("K113").Value = 3800.1426
$ws.Range("M113").Value = -1630.1426

Write-Host "Updated 168 cells across 7 sheets"